$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standard Excel PasteSpecial constants
$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Row 2: replace with the new production system record ---
# "1025" looks numeric, so stage it in a scratch cell formatted as Text,
# then paste-special just the value into A2. That keeps A2 storing a text
# string (matching the rest of the System ID column) instead of letting it
# collapse into a number.
$staging = $ws.Cells.Item(100, 100)
$staging.NumberFormat = "@"
$staging.Value = "1025"
$staging.Copy()
$ws.Cells.Item(2,1).PasteSpecial($xlPasteValues)
$staging.Clear()

$ws.Cells.Item(2,2).Value = "633 501D"
$ws.Cells.Item(2,3).Value = "172.17.85.63"

# --- Row 3: placeholder "bunk" entry across every column ---
$ws.Cells.Item(3,1).Value = "bunk"
$ws.Cells.Item(3,2).Value = "bunk"
$ws.Cells.Item(3,3).Value = "bunk"
$ws.Cells.Item(3,4).Value = "bunk"

# --- Row 4 (996 / 633_5thFl_TechLabCodecPro / ...) is no longer needed ---
$ws.Rows.Item(4).Delete()

# --- Drop the wrap-text formatting that used to be applied to the Status
# column. Pull clean formatting back from column A's header/data cells so
# the cells land back on the sheet's normal "Comma"/"Normal" styles instead
# of a brand-new derived one. ---
$ws.Cells.Item(1,1).Copy()
$ws.Cells.Item(1,4).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(2,4).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Cells.Item(2,4).Value = 200

# --- Let Excel recompute the (now-default) row heights for the data rows ---
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# --- Move the active selection to C3 ---
$ws.Range("C3").Select()
